$d = $word.ActiveDocument

# Locate the paragraph that ends with "...sufficient to lift" (the first
# half of the sentence that was split across two paragraphs) and the
# paragraph right after it that starts with "them and their families...".
$count = $d.Paragraphs.Count
$firstIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*sufficient to lift*") {
        $firstIdx = $i
        break
    }
}

if ($firstIdx -eq -1) {
    throw "Could not locate target paragraph."
}

$p1 = $d.Paragraphs.Item($firstIdx)
$p2 = $d.Paragraphs.Item($firstIdx + 1)

# Build the replacement: a single paragraph with three runs -
#   1) the original first-half text,
#   2) a standalone run containing just a space,
#   3) the original second-half text.
# Re-use the first paragraph's own pPr/attributes and each run's rsid so the
# result stays as close as possible to the original authoring metadata.
$p1Xml = $p1.Range.WordOpenXML
$rsidRPr = "000B7A87"
if ($p1Xml -match 'w:rsidRPr="([0-9A-Fa-f]+)"') {
    $rsidRPr = $matches[1]
}

$paraAttrsXml = ""
if ($p1Xml -match '(<w:p [^>]*>)') {
    $paraAttrsXml = $matches[1]
}

$text1 = $p1.Range.Text
$text1 = $text1.Substring(0, $text1.Length - 1)
$text2 = $p2.Range.Text
$text2 = $text2.Substring(0, $text2.Length - 1)

$newXml = '<w:p w14:paraId="7C0CB3C3" w14:textId="7ACCDCD5" w:rsidR="003A60DA" w:rsidRPr="' + $rsidRPr + '" w:rsidRDefault="003A60DA" w:rsidP="000B7A87" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="MText"/></w:pPr><w:r w:rsidRPr="' + $rsidRPr + '"><w:t>' + $text1 + '</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="' + $rsidRPr + '"><w:t>' + $text2 + '</w:t></w:r></w:p>'

$full = $d.Range($p1.Range.Start, $p2.Range.End)
$full.InsertXML($newXml)

# InsertXML only consumes the first paragraph of the target range; the old
# second paragraph survives as a leftover duplicate immediately after the
# newly-merged paragraph, so remove it explicitly.
$leftover = $d.Paragraphs.Item($firstIdx + 1)
$leftover.Range.Delete()
